$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.047044734400626
$ws.Range("D2").Value = 1.052466404352467
$ws.Range("E2").Value = 1.054340431559038
$ws.Range("F2").Value = 1.064609835100087
$ws.Range("I2").Value = 1.038029452620505
$ws.Range("J2").Value = 1.052095896572654
$ws.Range("K2").Value = 1.055215058576436
$ws.Range("L2").Value = 1.057083911176821
$ws.Range("M2").Value = 1.067325308703815
$ws.Range("N2").Value = 1.02112391404475
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.048272115889042
$ws.Range("D3").Value = 1.053424818691119
$ws.Range("E3").Value = 1.055421418262427
$ws.Range("F3").Value = 1.065751470073832
$ws.Range("I3").Value = 1.038258200611557
$ws.Range("J3").Value = 1.052970224662383
$ws.Range("K3").Value = 1.055985950634162
$ws.Range("L3").Value = 1.057977437212162
$ws.Range("M3").Value = 1.068281364930155
$ws.Range("N3").Value = 1.021422019144648
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.049066238001243
$ws.Range("D4").Value = 1.054044687540243
$ws.Range("E4").Value = 1.05612111705602
$ws.Range("F4").Value = 1.066490411094013
$ws.Range("I4").Value = 1.038404714845887
$ws.Range("J4").Value = 1.053535382884935
$ws.Range("K4").Value = 1.056483859536556
$ws.Range("L4").Value = 1.058555240284797
$ws.Range("M4").Value = 1.068899637731533
$ws.Range("N4").Value = 1.021614524083828
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.049400070545994
$ws.Range("D5").Value = 1.054305211671077
$ws.Range("E5").Value = 1.056415325838009
$ws.Range("F5").Value = 1.066801117040787
$ws.Range("I5").Value = 1.038465950209697
$ws.Range("J5").Value = 1.053772834888132
$ws.Range("K5").Value = 1.056692963213635
$ws.Range("L5").Value = 1.058798061353286
$ws.Range("M5").Value = 1.06915947429459
$ws.Range("N5").Value = 1.02169536021553
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.049456121588951
$ws.Range("D6").Value = 1.0543489508038
$ws.Range("E6").Value = 1.056464728059505
$ws.Range("F6").Value = 1.066853289194169
$ws.Range("I6").Value = 1.038476210835503
$ws.Range("J6").Value = 1.053812695909432
$ws.Range("K6").Value = 1.056728059919863
$ws.Range("L6").Value = 1.058838826979425
$ws.Range("M6").Value = 1.069203097015113
$ws.Range("N6").Value = 1.021708927512667
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.049070698750257
$ws.Range("D7").Value = 1.05404816894585
$ws.Range("E7").Value = 1.056125048069448
$ws.Range("F7").Value = 1.066494562544371
$ws.Range("I7").Value = 1.038405534487151
$ws.Range("J7").Value = 1.053538556281278
$ws.Range("K7").Value = 1.056486654442948
$ws.Range("L7").Value = 1.05855848521405
$ws.Range("M7").Value = 1.06890311001474
$ws.Range("N7").Value = 1.021615604585159
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047459550466824
$ws.Range("D8").Value = 1.052790365189931
$ws.Range("E8").Value = 1.054705709325667
$ws.Range("F8").Value = 1.064995609515127
$ws.Range("I8").Value = 1.038107069945116
$ws.Range("J8").Value = 1.052391502631021
$ws.Range("K8").Value = 1.055475773890704
$ws.Range("L8").Value = 1.057385959144166
$ws.Range("M8").Value = 1.067648487268583
$ws.Range("N8").Value = 1.021224740712406
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.044619842283873
$ws.Range("D9").Value = 1.050571714146237
$ws.Range("E9").Value = 1.052206363938162
$ws.Range("F9").Value = 1.0623559623771
$ws.Range("I9").Value = 1.03756963545884
$ws.Range("J9").Value = 1.050365686101175
$ws.Range("K9").Value = 1.053687478516552
$ws.Range("L9").Value = 1.055316965257354
$ws.Range("M9").Value = 1.065434894493468
$ws.Range("N9").Value = 1.020533001049828
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.042726148098697
$ws.Range("D10").Value = 1.049091070383678
$ws.Range("E10").Value = 1.05054123047221
$ws.Range("F10").Value = 1.060597284537548
$ws.Range("I10").Value = 1.037203603101668
$ws.Range("J10").Value = 1.049012018269528
$ws.Range("K10").Value = 1.052490539022368
$ws.Range("L10").Value = 1.053935666958425
$ws.Range("M10").Value = 1.063957242076702
$ws.Range("N10").Value = 1.020069819865231
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.041906000156953
$ws.Range("D11").Value = 1.048449560150633
$ws.Range("E11").Value = 1.049820456674344
$ws.Range("F11").Value = 1.059836002081074
$ws.Range("I11").Value = 1.037043268039823
$ws.Range("J11").Value = 1.048425110954639
$ws.Range("K11").Value = 1.051971116498695
$ws.Range("L11").Value = 1.053337070450381
$ws.Range("M11").Value = 1.063316933986045
$ws.Range("N11").Value = 1.019868774680948
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.041601333462485
$ws.Range("D12").Value = 1.04821121644196
$ws.Range("E12").Value = 1.049552764047153
$ws.Range("F12").Value = 1.059553262179297
$ws.Range("I12").Value = 1.036983435640133
$ws.Range("J12").Value = 1.048206991986013
$ws.Range("K12").Value = 1.051778007669247
$ws.Range("L12").Value = 1.053114651019609
$ws.Range("M12").Value = 1.063079022578139
$ws.Range("N12").Value = 1.019794024416156
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.041666686772137
$ws.Range("D13").Value = 1.048262344626705
$ws.Range("E13").Value = 1.049610183466511
$ws.Range("F13").Value = 1.05961390931614
$ws.Range("I13").Value = 1.03699628242418
$ws.Range("J13").Value = 1.048253784490065
$ws.Range("K13").Value = 1.051819437975047
$ws.Range("L13").Value = 1.053162364111044
$ws.Range("M13").Value = 1.063130058680559
$ws.Range("N13").Value = 1.019810061916146
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.041880816881855
$ws.Range("D14").Value = 1.048429859778704
$ws.Range("E14").Value = 1.049798328414226
$ws.Range("F14").Value = 1.059812630020188
$ws.Range("I14").Value = 1.037038327923666
$ws.Range("J14").Value = 1.048407083535218
$ws.Range("K14").Value = 1.051955157580367
$ws.Range("L14").Value = 1.053318686707722
$ws.Range("M14").Value = 1.063297269639854
$ws.Range("N14").Value = 1.019862597294436
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.042012745831924
$ws.Range("D15").Value = 1.048533063647628
$ws.Range("E15").Value = 1.049914255308559
$ws.Range("F15").Value = 1.05993507291183
$ws.Range("I15").Value = 1.037064196853883
$ws.Range("J15").Value = 1.048501520799545
$ws.Range("K15").Value = 1.052038756065494
$ws.Range("L15").Value = 1.053414992374443
$ws.Range("M15").Value = 1.063400284193144
$ws.Range("N15").Value = 1.01989495637476
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.042780574924624
$ws.Range("D16").Value = 1.049133637186394
$ws.Range("E16").Value = 1.050589070809106
$ws.Range("F16").Value = 1.060647813225885
$ws.Range("I16").Value = 1.037214205207369
$ws.Range("J16").Value = 1.049050953251981
$ws.Range("K16").Value = 1.052524987290409
$ws.Range("L16").Value = 1.05397538351722
$ws.Range("M16").Value = 1.063999727114816
$ws.Range("N16").Value = 1.020083152333611
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.043262167790081
$ws.Range("D17").Value = 1.049510258272776
$ws.Range("E17").Value = 1.051012428147108
$ws.Range("F17").Value = 1.061094958867567
$ws.Range("I17").Value = 1.037307808493452
$ws.Range("J17").Value = 1.049395393303989
$ws.Range("K17").Value = 1.05282968147179
$ws.Range("L17").Value = 1.05432677136547
$ws.Range("M17").Value = 1.06437561389399
$ws.Range("N17").Value = 1.020201072703433
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.043543056914816
$ws.Range("D18").Value = 1.049729898212838
$ws.Range("E18").Value = 1.051259388557
$ws.Range("F18").Value = 1.061355794275276
$ws.Range("I18").Value = 1.037362228131671
$ws.Range("J18").Value = 1.049596226094993
$ws.Range("K18").Value = 1.053007294520421
$ws.Range("L18").Value = 1.054531683212684
$ws.Range("M18").Value = 1.064594816555253
$ws.Range("N18").Value = 1.020269806917803
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04363883022923
$ws.Range("D19").Value = 1.049804783546527
$ws.Range("E19").Value = 1.051343599651956
$ws.Range("F19").Value = 1.061444736333753
$ws.Range("I19").Value = 1.037380753710572
$ws.Range("J19").Value = 1.049664692497159
$ws.Range("K19").Value = 1.053067837364733
$ws.Range("L19").Value = 1.054601544963179
$ws.Range("M19").Value = 1.064669551273927
$ws.Range("N19").Value = 1.02029323558968
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.043210499090284
$ws.Range("D20").Value = 1.049469854204934
$ws.Range("E20").Value = 1.050967003517451
$ws.Range("F20").Value = 1.061046982011326
$ws.Range("I20").Value = 1.037297784116509
$ws.Range("J20").Value = 1.049358445734064
$ws.Range("K20").Value = 1.052797002045995
$ws.Range("L20").Value = 1.054289075618417
$ws.Range("M20").Value = 1.064335289539889
$ws.Range("N20").Value = 1.020188425801089
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.041817761671919
$ws.Range("D21").Value = 1.048380532356717
$ws.Range("E21").Value = 1.049742923464043
$ws.Range("F21").Value = 1.059754110789126
$ws.Range("I21").Value = 1.037025954215865
$ws.Range("J21").Value = 1.048361943984586
$ws.Range("K21").Value = 1.051915196299744
$ws.Range("L21").Value = 1.053272655652105
$ws.Range("M21").Value = 1.063248032189707
$ws.Range("N21").Value = 1.019847128956274
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.040941933676985
$ws.Range("D22").Value = 1.047695295041141
$ws.Range("E22").Value = 1.0489734968065
$ws.Range("F22").Value = 1.058941429154278
$ws.Range("I22").Value = 1.036853442238767
$ws.Range("J22").Value = 1.047734735148788
$ws.Range("K22").Value = 1.051359773402232
$ws.Range("L22").Value = 1.052633163168347
$ws.Range("M22").Value = 1.062564010468122
$ws.Range("N22").Value = 1.019632118808455
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.0414062424849
$ws.Range("D23").Value = 1.048058584593077
$ws.Range("E23").Value = 1.0493813657802
$ws.Range("F23").Value = 1.05937222873535
$ws.Range("I23").Value = 1.036945046022539
$ws.Range("J23").Value = 1.048067294216126
$ws.Range("K23").Value = 1.051654308369891
$ws.Range("L23").Value = 1.052972211269371
$ws.Range("M23").Value = 1.062926663503638
$ws.Range("N23").Value = 1.019746139983286
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.043233846000673
$ws.Range("D24").Value = 1.04948811117777
$ws.Range("E24").Value = 1.050987528880645
$ws.Range("F24").Value = 1.061068660616843
$ws.Range("I24").Value = 1.03730231424969
$ws.Range("J24").Value = 1.049375140976983
$ws.Range("K24").Value = 1.052811768811556
$ws.Range("L24").Value = 1.054306108848435
$ws.Range("M24").Value = 1.064353510521205
$ws.Range("N24").Value = 1.020194140536353
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045354064571185
$ws.Range("D25").Value = 1.051145557872126
$ws.Range("E25").Value = 1.0528523073788
$ws.Range("F25").Value = 1.063038178752019
$ws.Range("I25").Value = 1.037709939156788
$ws.Range("J25").Value = 1.050889954604621
$ws.Range("K25").Value = 1.054150628720267
$ws.Range("L25").Value = 1.055852193094793
$ws.Range("M25").Value = 1.066007496710928
$ws.Range("N25").Value = 1.020712187680288
